# example.xlsx edit: add a "Cost" sheet (unit cost of IP links), remove the
# scratch "test" sheet, and refresh several saved cell-selection / active-tab
# view states left over from the author's last interactive session.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new "Cost" sheet right after "Links E(d,p)" (before "Beta")
#    BEFORE deleting "test", so the new sheet gets a fresh sheetId that
#    does not reuse the one still held by "test".
# ---------------------------------------------------------------------
$linksSheet = $wb.Worksheets.Item("Links E(d,p)")
$betaSheet  = $wb.Worksheets.Item("Beta")
$costSheet  = $wb.Worksheets.Add($betaSheet)
$costSheet.Name = "Cost"

# Populate the Cost sheet: header row + one row per IP link with a unit cost.
$costSheet.Range("A1").Value = "IP links E"
$costSheet.Range("B1").Value = "unit cost of IP link"
$costSheet.Range("A2").Value = "e1"
$costSheet.Range("B2").Value = 1
$costSheet.Range("A3").Value = "e2"
$costSheet.Range("B3").Value = 1
$costSheet.Range("A4").Value = "e3"
$costSheet.Range("B4").Value = 1

# Match the look of the workbook's other small lookup tables: centred values,
# the narrower 10pt font in the data column, and a wider column B.
$costSheet.Range("A1:A4").HorizontalAlignment = -4108
$costSheet.Range("B1:B4").HorizontalAlignment = -4108
$costSheet.Range("B1:B4").Font.Size = 10
$costSheet.Range("B1:B4").ColumnWidth = 17.85546875

# ---------------------------------------------------------------------
# 2. Remove the old scratch "test" sheet.
# ---------------------------------------------------------------------
$excel.DisplayAlerts = $false
$testSheet = $wb.Worksheets.Item("test")
$testSheet.Delete()
$excel.DisplayAlerts = $true

# ---------------------------------------------------------------------
# 3. Re-fetch worksheet handles by name (indices shifted after the
#    delete above) and refresh a few leftover cell selections.
# ---------------------------------------------------------------------
$costSheet = $wb.Worksheets.Item("Cost")
$costSheet.Activate()
$costSheet.Range("C24").Select()
$wb.Windows.Item(1).ScrollRow = 4

$linksSheet = $wb.Worksheets.Item("Links E(d,p)")
$linksSheet.Range("E6").Select()

$ipLinksSheet = $wb.Worksheets.Item("IPlinks")
$ipLinksSheet.Range("A1:B4").Select()

$pathsSheet = $wb.Worksheets.Item("Paths")
$pathsSheet.Range("H23").Select()

$pathsRdeSheet = $wb.Worksheets.Item("Paths R(d,e)")
$pathsRdeSheet.Range("H11").Select()

# ---------------------------------------------------------------------
# 4. Make "Paths R(d,e)" the active sheet/tab, matching the saved view.
# ---------------------------------------------------------------------
$pathsRdeSheet.Activate()
